{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph that separated it from the bibliography entry)\n// from the end of the document, mirroring the OOXML diff: the paragraph\n// \"Gest\u00e3o Empresarial - Estrat\u00e9gias Organizacionais Autor: Bertero, C. O.\n// Editora: ATLAS\" stays, the following blank paragraph + the \"Ver no\n// Jupiter...\" paragraph + the \"\u00a9 2020...\" paragraph are deleted, while the\n// blank paragraph right before the final page-break paragraph is kept.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the bibliography entry that must remain immediately before the\n// block we delete.\nconst anchorText =\n  \"Gest\u00e3o Empresarial - Estrat\u00e9gias Organizacionais Autor: Bertero, C. O. Editora: ATLAS\";\nconst anchorIndex = items.findIndex((p) => p.text === anchorText);\nif (anchorIndex === -1) {\n  throw new Error(\"Could not locate the anchor bibliography paragraph.\");\n}\n\n// Collect the paragraphs to delete: every paragraph after the anchor up to\n// (and not including) the first one that follows the \"\u00a9 2020\" copyright\n// paragraph.\nconst footerText =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\nconst footerIndex = items.findIndex((p) => p.text === footerText);\nif (footerIndex === -1) {\n  throw new Error(\"Could not locate the footer paragraph.\");\n}\n\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i <= footerIndex; i++) {\n  toDelete.push(items[i]);\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph that separated it from the bibliography entry)\n# from the end of the document, mirroring the OOXML diff: the paragraph\n# \"Gest\u00e3o Empresarial - Estrat\u00e9gias Organizacionais Autor: Bertero, C. O.\n# Editora: ATLAS\" stays, the following blank paragraph + the \"Ver no\n# Jupiter...\" paragraph + the \"\u00a9 2020...\" paragraph are deleted, while the\n# blank paragraph right before the final page-break paragraph is kept.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Gest\u00e3o Empresarial - Estrat\u00e9gias Organizacionais Autor: Bertero, C. O. Editora: ATLAS`r\"\n$footerText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution`r\"\n\n$anchorIndex = -1\n$footerIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -eq $anchorText) { $anchorIndex = $i }\n    if ($t -eq $footerText) { $footerIndex = $i }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not locate the anchor bibliography paragraph.\"\n}\nif ($footerIndex -eq -1) {\n    throw \"Could not locate the footer paragraph.\"\n}\n\n# Delete from the end backwards so earlier indices stay valid.\nfor ($i = $footerIndex; $i -gt $anchorIndex; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n}\n"}
